# Sprint 3 config sheet update: add CertificateTemplatesPath parameter row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above the existing "CountryLookUpPath" row (row 26) so every
# row below shifts down by one (A1:C34 -> A1:C35).
$ws.Rows.Item(26).Insert()

# Populate the new row. Set Value/Description before Name so the shared
# string table gets the three new entries in the same order as the source
# edit (path, description, name).
$ws.Range("B26").Value2 = "\\EARTH.GSI.GOV.UK\USER\SHARED\Agency\CoFS for G drive\RobotDocuments\Robot Certificate Templates\"
$ws.Range("C26").Value2 = "Folder path for the robot's version of certificate templates"
$ws.Range("A26").Value2 = "CertificateTemplatesPath"

# Match the row height used by the other wrapped-text rows in this block.
$ws.Rows.Item(26).RowHeight = 48.75

# Grow Table1 so the new row is included in the table / autofilter range.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:C35"))

# Sheet1 becomes the active / selected tab (was ApplicationFields before),
# with the whole new row selected.
$ws.Activate()
$ws.Rows.Item(26).Select()
